$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the time value that was stored in K2 ("08:14:14 AM"), leaving the
# cell's formatting/style untouched.
$ws.Range("K2").ClearContents()

# Move/update the active selection on the sheet to K2.
$ws.Range("K2").Select()
